$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update DataHora (column C) timestamps for rows 2-11.
# Rows 2-5 previously "24/07/2025 13:29:39" -> "25/07/2025 10:45:12"
$ws.Range("C2").Value = "25/07/2025 10:45:12"
$ws.Range("C3").Value = "25/07/2025 10:45:12"
$ws.Range("C4").Value = "25/07/2025 10:45:12"
$ws.Range("C5").Value = "25/07/2025 10:45:12"

# Rows 6-11 previously "24/07/2025 13:29:53" -> "25/07/2025 10:45:29"
$ws.Range("C6").Value = "25/07/2025 10:45:29"
$ws.Range("C7").Value = "25/07/2025 10:45:29"
$ws.Range("C8").Value = "25/07/2025 10:45:29"
$ws.Range("C9").Value = "25/07/2025 10:45:29"
$ws.Range("C10").Value = "25/07/2025 10:45:29"
$ws.Range("C11").Value = "25/07/2025 10:45:29"
